$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C:G block for rows 1-15
$cg = New-Object 'object[,]' 15,5
$cg[0,0] = 306.50408290847901
$cg[0,1] = 306.50408290847901
$cg[0,2] = 306.50408290847901
$cg[0,3] = 306.50408290848014
$cg[0,4] = 2000
$cg[1,0] = 748.71154169204465
$cg[1,1] = 2129.795400180607
$cg[1,2] = 1914.6285076969839
$cg[1,3] = 911.81207739129502
$cg[1,4] = 5000
$cg[2,0] = 1423.7622406159676
$cg[2,1] = 5129.7954001806065
$cg[2,2] = 2585.5359637645511
$cg[2,3] = 1604.0879345672024
$cg[2,4] = 8000
$cg[3,0] = 2140.223972604711
$cg[3,1] = 5869.7675198705365
$cg[3,2] = 3325.5080834544806
$cg[3,3] = 2354.4243853174103
$cg[3,4] = 11000
$cg[4,0] = 2717.1177608244061
$cg[4,1] = 6433.7129346635647
$cg[4,2] = 3889.4534982475088
$cg[4,3] = 2960.9831667057238
$cg[4,4] = 14000
$cg[5,0] = 3207.5564041495486
$cg[5,1] = 6924.1515779887068
$cg[5,2] = 4379.8921415726509
$cg[5,3] = 3503.3690674977679
$cg[5,4] = 17000
$cg[6,0] = 3826.5315741711356
$cg[6,1] = 7543.1267480102942
$cg[6,2] = 4995.7402857209872
$cg[6,3] = 4168.8188419005755
$cg[6,4] = 20000
$cg[7,0] = 4447.7164724006834
$cg[7,1] = 8184.9283287885537
$cg[7,2] = 5619.3862321498691
$cg[7,3] = 4868.3692482874394
$cg[7,4] = 23000
$cg[8,0] = 4743.8792325366303
$cg[8,1] = 8481.0910889245006
$cg[8,2] = 5915.548992285816
$cg[8,3] = 5283.543671113066
$cg[8,4] = 26000
$cg[9,0] = 5324.0519372706467
$cg[9,1] = 11481.091088924501
$cg[9,2] = 6485.1725979249213
$cg[9,3] = 5957.3234459801306
$cg[9,4] = 29000
$cg[10,0] = 5811.408529753634
$cg[10,1] = 14481.091088924501
$cg[10,2] = 7520.7655216721096
$cg[10,3] = 6878.6636113148788
$cg[10,4] = 32000
$cg[11,0] = 6344.5945747687656
$cg[11,1] = 17722.779321371192
$cg[11,2] = 8480.2491672083997
$cg[11,3] = 8076.8521740019578
$cg[11,4] = 35000
$cg[12,0] = 7067.6327411904213
$cg[12,1] = 18445.81748779285
$cg[12,2] = 9203.2873336300545
$cg[12,3] = 9101.529322978311
$cg[12,4] = 38000
$cg[13,0] = 8240.6378593561785
$cg[13,1] = 21843.723122207048
$cg[13,2] = 9995.7996752080726
$cg[13,3] = 9989.5130500311716
$cg[13,4] = 41000
$cg[14,0] = 8513.208248724226
$cg[14,1] = 22262.922557139158
$cg[14,2] = 10413.075965441163
$cg[14,3] = 10413.075965441165
$cg[14,4] = 44000
$ws.Range("C1:G15").Value = $cg

# Update I:L block for rows 1-15
$il = New-Object 'object[,]' 15,4
$il[0,0] = 15.044217687074831
$il[0,1] = 15.044217687074831
$il[0,2] = 15.044217687074831
$il[0,3] = 15.044217687074729
$il[1,0] = 20.044217687074831
$il[1,1] = 30
$il[1,2] = 30
$il[1,3] = 20.049259636141919
$il[2,0] = 25.044217687074831
$il[2,1] = 45
$il[2,2] = 35
$il[2,3] = 25.054231228411776
$il[3,0] = 30.044217687074831
$il[3,1] = 50
$il[3,2] = 40
$il[3,3] = 30.066287547394531
$il[4,0] = 35.044217687074834
$il[4,1] = 55
$il[4,2] = 45
$il[4,3] = 35.069828928961428
$il[5,0] = 40.044217687074834
$il[5,1] = 60
$il[5,2] = 50
$il[5,3] = 40.072527122558661
$il[6,0] = 45.044217687074827
$il[6,1] = 65
$il[6,2] = 55
$il[6,3] = 45.074852003245766
$il[7,0] = 50.044217687074827
$il[7,1] = 70
$il[7,2] = 60
$il[7,3] = 50.08073183755355
$il[8,0] = 55.052403846153844
$il[8,1] = 75
$il[8,2] = 65
$il[8,3] = 55.081660485062379
$il[9,0] = 60.015890724937918
$il[9,1] = 90
$il[9,2] = 70
$il[9,3] = 60.03480390112599
$il[10,0] = 65.012284935195325
$il[10,1] = 105
$il[10,2] = 75
$il[10,3] = 65.966991355439802
$il[11,0] = 70.012284935195325
$il[11,1] = 120
$il[11,2] = 80
$il[11,3] = 75.649238649397986
$il[12,0] = 75.012284935195325
$il[12,1] = 125
$il[12,2] = 85
$il[12,3] = 83.718637537139458
$il[13,0] = 80.055576366532222
$il[13,1] = 140
$il[13,2] = 90
$il[13,3] = 89.910630282508734
$il[14,0] = 85.114016693392543
$il[14,1] = 145
$il[14,2] = 95
$il[14,3] = 95
$ws.Range("I1:L15").Value = $il

# Update O:R block for rows 1-15
$orr = New-Object 'object[,]' 15,4
$orr[0,0] = 18.481891713230183
$orr[0,1] = 0.005359497648784229
$orr[0,2] = 0.0056564573953785781
$orr[0,3] = 0.0025075241180238638
$orr[1,0] = 24.795251659163029
$orr[1,1] = 0.0095018585584140305
$orr[1,2] = 0.0093371653656188873
$orr[1,3] = 0.0029640269106970363
$orr[2,0] = 33.439194931886995
$orr[2,1] = 0.012132398535908804
$orr[2,2] = 0.012734566910947346
$orr[2,3] = 0.0052973019329690179
$orr[3,0] = 51.454941194005734
$orr[3,1] = 0.029356348282582798
$orr[3,2] = 0.030817822591013611
$orr[3,3] = 0.0077927192524351267
$orr[4,0] = 61.567998946418676
$orr[4,1] = 0.020373884836507161
$orr[4,2] = 0.020351698188773103
$orr[4,3] = 0.010642245140839701
$orr[5,0] = 36.736340704989267
$orr[5,1] = 0.02640979079697852
$orr[5,2] = 0.027914216179868862
$orr[5,3] = 0.015177203902119336
$orr[6,0] = 71.583105209083556
$orr[6,1] = 0.02683617265535489
$orr[6,2] = 0.027542732052424204
$orr[6,3] = 0.021834861652029167
$orr[7,0] = 173.50511572896792
$orr[7,1] = 0.039251593171973825
$orr[7,2] = 0.033056967346943418
$orr[7,3] = 0.047227585228016157
$orr[8,0] = 198.47734032155847
$orr[8,1] = 0.037620305675116934
$orr[8,2] = 0.043703429373073602
$orr[8,3] = 0.12199893617868557
$orr[9,0] = 1125.1184610544688
$orr[9,1] = 0.03804042976105549
$orr[9,2] = 0.035795027232687876
$orr[9,3] = 0.098273851784090926
$orr[10,0] = 2354.4699382300973
$orr[10,1] = 0.038662793674416066
$orr[10,2] = 0.038940980103696982
$orr[10,3] = 0.30359976396820137
$orr[11,0] = 3685.2387143918527
$orr[11,1] = 0.047465772829207185
$orr[11,2] = 0.048170625561066188
$orr[11,3] = 0.54541890410475746
$orr[12,0] = 306.63783796237828
$orr[12,1] = 0.046577738031441322
$orr[12,2] = 0.047711248175290447
$orr[12,3] = 0.61472545246539168
$orr[13,0] = 7200.8351617406615
$orr[13,1] = 0.052182995470510529
$orr[13,2] = 0.050872845477394082
$orr[13,3] = 0.38322450014620429
$orr[14,0] = 7200.219314728407
$orr[14,1] = 0.056167205403984719
$orr[14,2] = 0.055471739330782435
$orr[14,3] = 0.37917724287097498
$ws.Range("O1:R15").Value = $orr

# Update A column for rows 11-20
$acol = New-Object 'object[,]' 10,1
$acol[0,0] = 11
$acol[1,0] = 12
$acol[2,0] = 13
$acol[3,0] = 14
$acol[4,0] = 15
$acol[5,0] = 16
$acol[6,0] = 17
$acol[7,0] = 18
$acol[8,0] = 19
$acol[9,0] = 20
$ws.Range("A11:A20").Value = $acol

# Update M column for rows 11-15
$mcol = New-Object 'object[,]' 5,1
$mcol[0,0] = 70
$mcol[1,0] = 75
$mcol[2,0] = 80
$mcol[3,0] = 85
$mcol[4,0] = 90
$ws.Range("M11:M15").Value = $mcol
